# Generate Report for Handoff
# - Status moves from "Handed back: in sync with en-US" to "Ready for handoff"
#   on every sheet that tracks it.
# - The "Latest Handoff Datetime" / "Latest HO Xliff Generate Date" stamps
#   advance to the new handoff-generation timestamps.
# - The now-shorter Status column/cells autosize narrower.

$wb = $excel.ActiveWorkbook

# Target narrower Status-column width (matches the handed-off report's
# autofit to the shorter "Ready for handoff" text).
$statusColWidth = 16.3333333333333

# --- Overview sheet --------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-26 02:59:38"

$wsOverview.Columns.Item(5).ColumnWidth = $statusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $statusColWidth

# --- zh-cn sheet -------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Ready for handoff"
$wsZhCn.Range("H2").Value = "2016-08-26 02:59:34"

$wsZhCn.Columns.Item(3).ColumnWidth = $statusColWidth

# --- de-de sheet -------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Ready for handoff"
$wsDeDe.Range("H2").Value = "2016-08-26 02:59:38"

$wsDeDe.Columns.Item(3).ColumnWidth = $statusColWidth
